# "add WishListTest with data in xlsx"
#
# Inserts a new worksheet "addToWishList" right after "signIn", populates it
# with the same login form layout used elsewhere in the workbook (url /
# login / password headers, homepage url, test e-mail w/ mailto hyperlink,
# and a password), and makes it the active sheet. Also adds a hyperlink
# (to the shop's homepage) on the "hoverToJackets" sheet's existing url cell.

$wb = $excel.ActiveWorkbook

# --- insert the new sheet right after "signIn" --------------------------
$signIn = $wb.Worksheets.Item("signIn")
$wishList = $wb.Worksheets.Add($null, $signIn)
$wishList.Name = "addToWishList"

# --- header row -----------------------------------------------------------
$wishList.Range("A1").Value = "url"
$wishList.Range("B1").Value = "login"
$wishList.Range("C1").Value = "password"

# --- data row ---------------------------------------------------------
$wishList.Range("A2").Value = "https://magento.softwaretestingboard.com/"
$wishList.Range("B2").Value = "sistulostu@gufum.com"
$wishList.Range("C2").Value = "password!123"

# e-mail cell gets a mailto hyperlink + vertical-center alignment, matching
# the styling used for the same data on the "signIn" sheet
$wishList.Hyperlinks.Add($wishList.Range("B2"), "mailto:sistulostu@gufum.com") | Out-Null
$wishList.Range("B2").VerticalAlignment = -4108

$wishList.PageSetup.Orientation = 1

$wishList.Range("C2").Select() | Out-Null

# --- hoverToJackets: turn the url cell into a real hyperlink ------------
$hover = $wb.Worksheets.Item("hoverToJackets")
$hover.Hyperlinks.Add($hover.Range("A2"), "https://magento.softwaretestingboard.com/") | Out-Null
$hover.Range("A2").VerticalAlignment = -4108

Write-Host "addToWishList sheet added"
